$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 177: odds refreshed for existing match (id 8240846) ---
$ws.Range("M177").Value = 2.75
$ws.Range("O177").Value = 2.55
$ws.Range("P177").Value = 0
$ws.Range("Q177").Value = 2
$ws.Range("R177").Value = 1.85
$ws.Range("S177").Value = 2.5
$ws.Range("T177").Value = 2.025
$ws.Range("U177").Value = 1.825

# --- Row 179: now holds the match that used to live in row 181 (id 7129721) ---
$ws.Range("B179").Value = "'7129721"
$ws.Range("B179").Style = "Normal"
$ws.Range("E179").Value = "NK Lokomotiva Zagreb"
$ws.Range("F179").Value = "Hajduk Split"
$ws.Range("J179").Value = 3.3
$ws.Range("K179").Value = 3.4
$ws.Range("L179").Value = 2.05
$ws.Range("M179").Value = 3.5
$ws.Range("N179").Value = 3.5
$ws.Range("O179").Value = 1.95
$ws.Range("P179").Value = 0.5
$ws.Range("Q179").Value = 1.8
$ws.Range("R179").Value = 2.05
$ws.Range("S179").Value = 2.75
$ws.Range("T179").Value = 2.025
$ws.Range("U179").Value = 1.825

# --- Row 180: now holds the match that used to live in row 179 (id 8240844) ---
$ws.Range("B180").Value = "'8240844"
$ws.Range("B180").Style = "Normal"
$ws.Range("D180").Value = 45438.41666666666
$ws.Range("E180").Value = "HNK Rijeka"
$ws.Range("F180").Value = "Slaven Belupo"
$ws.Range("J180").Value = 1.42
$ws.Range("K180").Value = 5
$ws.Range("L180").Value = 6
$ws.Range("M180").Value = 1.42
$ws.Range("N180").Value = 5
$ws.Range("O180").Value = 6
$ws.Range("P180").Value = -1.25
$ws.Range("Q180").Value = 1.9
$ws.Range("R180").Value = 1.95
$ws.Range("S180").Value = 3

# --- Row 181: now holds the match that used to live in row 180 (id 8240845) ---
$ws.Range("B181").Value = "'8240845"
$ws.Range("B181").Style = "Normal"
$ws.Range("D181").Value = 45438.52083333334
$ws.Range("E181").Value = "Dinamo Zagreb"
$ws.Range("F181").Value = "NK Rudes"
$ws.Range("J181").Value = 1.071
$ws.Range("K181").Value = 11
$ws.Range("L181").Value = 21
$ws.Range("M181").Value = 1.062
$ws.Range("N181").Value = 11
$ws.Range("O181").Value = 23
$ws.Range("P181").Value = -3
$ws.Range("Q181").Value = 1.925
$ws.Range("R181").Value = 1.925
$ws.Range("S181").Value = 4.25
$ws.Range("T181").Value = 1.925
$ws.Range("U181").Value = 1.925
